$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Summary block (rows 10-12): the grading logic now produces real numbers
# instead of the previous all-zero / "Absent" placeholder values.
# ---------------------------------------------------------------------------

# A10/A11/A12 labels pick up the same bold header style ("s=4") already used
# by the row-9 header cells.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 10: No. Right / Wrong / Not Attempt / Max
$ws.Range("B10").Value = 17
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 7
$ws.Range("E10").Value = 28

# Row 11: Marking scheme (+4 for right, -1 for wrong). C11 used to hold the
# "-1" as text; it must become a real number now.
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12: Totals + final score fraction
$ws.Range("B12").Value = 68
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "64/112"

# ---------------------------------------------------------------------------
# Answer-key grid: the sheet used to have three Student-Ans/Correct-Ans
# column pairs (A:B, D:E, G:H). The third pair (G:H) is removed entirely,
# and the second pair (D:E) is trimmed down to only the first three
# questions. The remaining cells in the first pair (A:B) get colored to
# show whether the student's (now-populated) answer was right (green),
# wrong (red), or left blank (plain/black - unchanged).
# ---------------------------------------------------------------------------

# Drop the whole third column-pair (Student Ans / Correct Ans) for rows 15-40.
$ws.Range("G15:H40").Clear()

# Drop the second column-pair's rows beyond the first three questions.
$ws.Range("D19:E40").Clear()

# Reference cells already present on the sheet that carry the colors we need:
#   B10 -> green "correct" style (s=5)
#   C10 -> red "incorrect" style (s=6)

# D16:D18 (kept Student-Ans column of the second pair) are now all correct -> green.
$ws.Range("B10").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option C"
$ws.Range("D18").Value = "Option D"

# Column A (first pair) - correctly answered questions (green / style s=5,
# copied from the already-green-styled B10 cell).
$ws.Range("B10").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A16").Value = "Option A"
$ws.Range("A18").Value = "Option B"
$ws.Range("A19").Value = "Option C"
$ws.Range("A22").Value = "Option D"
$ws.Range("A23").Value = "Option D"
$ws.Range("A25").Value = "Option A"
$ws.Range("A27").Value = "Option A"
$ws.Range("A28").Value = "Option D"
$ws.Range("A30").Value = "Option B"
$ws.Range("A32").Value = "Option C"
$ws.Range("A33").Value = "Option D"
$ws.Range("A35").Value = "Option D"
$ws.Range("A36").Value = "Option A"
$ws.Range("A39").Value = "Option D"

# Column A (first pair) - incorrectly answered questions (red / style s=6,
# copied from the already-red-styled C10 cell).
$ws.Range("C10").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("A34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A20").Value = "Option A"
$ws.Range("A21").Value = "Option D"
$ws.Range("A26").Value = "Option D"
$ws.Range("A34").Value = "Option A"

# Rows A17, A24, A29, A31, A37, A38, A40 remain blank / not attempted, no
# change needed there (they already carry the plain "not attempted" style).
